# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.746.13"
$ws.Range("E2").Value = "  +1.86%  "

$ws.Range("D3").Value = "3.773.15"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.06%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.69"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.37%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.39"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "3.765.77"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  -0.02%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("E10").Value = "  -2.40%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.50"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.16%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("E13").Value = "  -4.90%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.65"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").Value = "4.427.02"
$ws.Range("E15").Value = "  +0.61%  "

$ws.Range("D16").Value = "3.783.63"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").Value = "68.870.01"
$ws.Range("E17").Value = "  +1.93%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.09"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -2.63%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("E20").Value = "  -0.23%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +3.76%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.41"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.54%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.706"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -1.93%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.59"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.75%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.25%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +0.40%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.05%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("D30").Value = "3.938.18"
$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("E31").Value = "  -3.17%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.42"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -3.11%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.19"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.23%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.22"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -0.48%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.34"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("D37").Value = "3.745.29"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  -2.18%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.50"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -8.27%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.139"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +1.50%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +0.31%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.85"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.308"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.15%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +1.33%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.73"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +12.09%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.63"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -1.32%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.02"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.19%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "398.13"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.17%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "146.26"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +5.35%  "
